# T2244 Contacts - Duplicate Management test data update
# Update the Email address on the "Contact" sheet (row 3) from
# test@gmail.com to testtest@mailinator.com, and update the active
# cell selection as it was left after editing, matching the commit
# "Contacts & Reports changes latest - 30 Mar 2023".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Contact")

$ws.Range("E3").Value = "testtest@mailinator.com"

$ws.Activate()
$ws.Range("D13").Select() | Out-Null
